$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1770.6487
$ws.Range("I17").Value = 1416.5
$ws.Range("J17").Value = 1839.1936
$ws.Range("K17").Value = 4249.5
$ws.Range("L17").Value = 5517.5808
$ws.Range("M17").Value = -4081.5
$ws.Range("N17").Value = -5853.5808

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 355.83334
$ws.Range("I19").Value = 324.33334
$ws.Range("J19").Value = 387.33334
$ws.Range("K19").Value = 324.33334
$ws.Range("L19").Value = 387.33334
$ws.Range("M19").Value = -149.33334
$ws.Range("N19").Value = -737.33334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 287.46155
$ws.Range("I42").Value = 206.25
$ws.Range("K42").Value = 618.75
$ws.Range("M42").Value = -388.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9263669
$ws.Range("I51").Value = 4894
$ws.Range("J51").Value = 10421016
$ws.Range("K51").Value = 4894
$ws.Range("L51").Value = 10421016
$ws.Range("M51").Value = -4410
$ws.Range("N51").Value = -10421984

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 84153.25
$ws.Range("J57").Value = 84153.25
$ws.Range("L57").Value = 252459.75
$ws.Range("N57").Value = -253457.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6185.6
$ws.Range("I62").Value = 5540.1113
$ws.Range("J62").Value = 11995
$ws.Range("K62").Value = 5540.1113
$ws.Range("L62").Value = 11995
$ws.Range("M62").Value = -4916.1113
$ws.Range("N62").Value = -13243

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 6185.6
$ws.Range("I65").Value = 5540.1113
$ws.Range("J65").Value = 11995
$ws.Range("K65").Value = 27700.5565
$ws.Range("L65").Value = 59975
$ws.Range("M65").Value = -24580.5565
$ws.Range("N65").Value = -66215

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2225.4119
$ws.Range("J80").Value = 2069.1667
$ws.Range("L80").Value = 6207.500100000001
$ws.Range("N80").Value = -8203.500100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2225.4119
$ws.Range("J83").Value = 2069.1667
$ws.Range("L83").Value = 18622.5003
$ws.Range("N83").Value = -28606.5003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1060.4286
$ws.Range("I98").Value = 1060.4286
$ws.Range("K98").Value = 1060.4286
$ws.Range("M98").Value = 437.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 2202998.2
$ws.Range("J101").Value = 336664.34
$ws.Range("L101").Value = 1009993.02
$ws.Range("N101").Value = -1013237.02

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 8087.684
$ws.Range("I106").Value = 2804.1765
$ws.Range("K106").Value = 2804.1765
$ws.Range("M106").Value = -2173.1765

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1060.4286
$ws.Range("I122").Value = 1060.4286
$ws.Range("K122").Value = 3181.2858
$ws.Range("M122").Value = -731.2857999999997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 78343.75
$ws.Range("J133").Value = 78343.75
$ws.Range("L133").Value = 78343.75
$ws.Range("N133").Value = -88463.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 112968.62
$ws.Range("J134").Value = 107424.5
$ws.Range("L134").Value = 107424.5
$ws.Range("N134").Value = -117564.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3388.3914
$ws.Range("I137").Value = 2682.1667
$ws.Range("K137").Value = 8046.500100000001
$ws.Range("M137").Value = -5496.500100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5270.183
$ws.Range("J138").Value = 5330.418
$ws.Range("L138").Value = 15991.254
$ws.Range("N138").Value = -26271.254

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14162.652
$ws.Range("I32").Value = 11822.19
$ws.Range("J32").Value = 38737.5
$ws.Range("K32").Value = 11822.19
$ws.Range("L32").Value = 38737.5
$ws.Range("M32").Value = -11535.19
$ws.Range("N32").Value = -39311.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10735.071
$ws.Range("I45").Value = 17072.428
$ws.Range("J45").Value = 4397.7144
$ws.Range("K45").Value = 17072.428
$ws.Range("L45").Value = 4397.7144
$ws.Range("M45").Value = -16695.428
$ws.Range("N45").Value = -5151.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2041.5625
$ws.Range("I74").Value = 1840.3572
$ws.Range("K74").Value = 1840.3572
$ws.Range("M74").Value = -966.3571999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2041.5625
$ws.Range("I77").Value = 1840.3572
$ws.Range("K77").Value = 9201.786
$ws.Range("M77").Value = -4833.786

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4043.95
$ws.Range("I122").Value = 4045.5264
$ws.Range("K122").Value = 12136.5792
$ws.Range("M122").Value = -9686.5792

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 119843
$ws.Range("J58").Value = 119843
$ws.Range("L58").Value = 119843
$ws.Range("N58").Value = -120431

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2697
$ws.Range("I105").Value = 2697
$ws.Range("K105").Value = 2697
$ws.Range("M105").Value = -950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 129012.91
$ws.Range("J132").Value = 129012.91
$ws.Range("L132").Value = 129012.91
$ws.Range("N132").Value = -139132.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 36352.484
$ws.Range("J31").Value = 4928.2104
$ws.Range("L31").Value = 4928.2104
$ws.Range("N31").Value = -5518.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 36352.484
$ws.Range("J34").Value = 4928.2104
$ws.Range("L34").Value = 4928.2104
$ws.Range("N34").Value = -5332.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2514.4443
$ws.Range("J94").Value = 2181
$ws.Range("L94").Value = 2181
$ws.Range("N94").Value = -3083

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4319.3335
$ws.Range("I134").Value = 4279.294
$ws.Range("K134").Value = 12837.882
$ws.Range("M134").Value = -10302.882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 881284.4399999999
$ws.Range("J141").Value = 881284.4399999999
$ws.Range("L141").Value = 881284.4399999999
$ws.Range("N141").Value = -891644.4399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9010.85
$ws.Range("I56").Value = 9010.85
$ws.Range("K56").Value = 9010.85
$ws.Range("M56").Value = -8480.85

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 5001
$ws.Range("J80").Value = 5001
$ws.Range("L80").Value = 15003
$ws.Range("N80").Value = -16875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 5001
$ws.Range("J83").Value = 5001
$ws.Range("L83").Value = 45009
$ws.Range("N83").Value = -54369

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2195.8125
$ws.Range("I102").Value = 2096.3845
$ws.Range("K102").Value = 2096.3845
$ws.Range("M102").Value = -474.3845000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4146.0713
$ws.Range("I122").Value = 4105
$ws.Range("J122").Value = 4248.75
$ws.Range("K122").Value = 12315
$ws.Range("L122").Value = 12746.25
$ws.Range("M122").Value = -9865
$ws.Range("N122").Value = -17646.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 734.6875
$ws.Range("J16").Value = 529
$ws.Range("L16").Value = 529
$ws.Range("N16").Value = -869

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 405081.6
$ws.Range("I122").Value = 670869.3
$ws.Range("K122").Value = 2012607.9
$ws.Range("M122").Value = -2010157.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 31384.85
$ws.Range("I136").Value = 38512.8
$ws.Range("J136").Value = 10001
$ws.Range("K136").Value = 115538.4
$ws.Range("L136").Value = 30003
$ws.Range("M136").Value = -112988.4
$ws.Range("N136").Value = -35103

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 100599.4
$ws.Range("J139").Value = 100599.4
$ws.Range("L139").Value = 100599.4
$ws.Range("N139").Value = -110879.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 29412946
$ws.Range("I107").Value = 1317.6923
$ws.Range("J107").Value = 125000744
$ws.Range("K107").Value = 3953.0769
$ws.Range("L107").Value = 375002232
$ws.Range("M107").Value = -2033.0769
$ws.Range("N107").Value = -375006072

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 615.6667
$ws.Range("I113").Value = 649.5
$ws.Range("J113").Value = 548
$ws.Range("K113").Value = 1948.5
$ws.Range("L113").Value = 1644
$ws.Range("M113").Value = 221.5
$ws.Range("N113").Value = -5984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 20000
$ws.Range("J119").Value = 20000
$ws.Range("L119").Value = 20000
$ws.Range("N119").Value = -29676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1829.1666
$ws.Range("I122").Value = 1695
$ws.Range("K122").Value = 5085
$ws.Range("M122").Value = -2635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3114.6667
$ws.Range("I126").Value = 3114.6667
$ws.Range("K126").Value = 9344.000100000001
$ws.Range("M126").Value = -6874.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 86472.42999999999
$ws.Range("J140").Value = 86472.42999999999
$ws.Range("L140").Value = 86472.42999999999
$ws.Range("N140").Value = -96832.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 111111
$ws.Range("J141").Value = 111111
$ws.Range("L141").Value = 111111
$ws.Range("N141").Value = -121471
